$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '29.423.70'
$ws.Cells.Item(2, 5).Value = '  +0.07%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.848.73'
$ws.Cells.Item(3, 5).Value = '  +0.28%  '

$ws.Cells.Item(4, 5).Value = '  +0.11%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '240.50'
$ws.Cells.Item(5, 5).Value = '  +0.57%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.6265'
$ws.Cells.Item(6, 5).Value = '  -0.87%  '

$ws.Cells.Item(7, 5).Value = '  +0.06%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.07695'
$ws.Cells.Item(8, 5).Value = '  +2.21%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.2913'
$ws.Cells.Item(9, 5).Value = '  -0.46%  '

$ws.Cells.Item(10, 5).Value = '  +1.45%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07746'
$ws.Cells.Item(11, 5).Value = '  +0.40%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '1.845.77'
$ws.Cells.Item(12, 5).Value = '  -0.06%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '5.024'
$ws.Cells.Item(13, 5).Value = '  +0.46%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.00001076'
$ws.Cells.Item(15, 5).Value = '  +3.51%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '83.52'
$ws.Cells.Item(16, 5).Value = '  +0.43%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '6.163'
$ws.Cells.Item(17, 5).Value = '  -0.07%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '29.446.21'
$ws.Cells.Item(18, 5).Value = '  +0.06%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '228.40'
$ws.Cells.Item(19, 5).Value = '  +0.05%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '12.37'
$ws.Cells.Item(20, 5).Value = '  -0.26%  '

$ws.Cells.Item(21, 5).Value = '  +0.06%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '7.409'
$ws.Cells.Item(22, 5).Value = '  -0.54%  '

$ws.Cells.Item(23, 5).Value = '  +0.02%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '157.04'
$ws.Cells.Item(24, 5).Value = '  -0.09%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.1374'
$ws.Cells.Item(25, 5).Value = '  -1.40%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '8.395'

$ws.Cells.Item(27, 5).Value = '  +0.61%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '1.341'
$ws.Cells.Item(28, 5).Value = '  +4.66%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '1.463'
$ws.Cells.Item(29, 5).Value = '  +0.32%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.05643'
$ws.Cells.Item(30, 5).Value = '  +0.21%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '4.115'
$ws.Cells.Item(31, 5).Value = '  +0.38%  '

$ws.Cells.Item(32, 5).Value = '  +0.06%  '

$ws.Cells.Item(33, 5).Value = '  +0.07%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.160'
$ws.Cells.Item(34, 5).Value = '  +0.25%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.7077'
$ws.Cells.Item(35, 5).Value = '  -0.27%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.597'
$ws.Cells.Item(36, 5).Value = '  +0.33%  '

$ws.Cells.Item(37, 2).Value = 'Maker'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.228.27'
$ws.Cells.Item(37, 5).Value = '  -1.40%  '

$ws.Cells.Item(38, 2).Value = 'MXToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.767'
$ws.Cells.Item(38, 5).Value = '  +0.06%  '

$ws.Cells.Item(39, 5).Value = '  -1.12%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '6.464'
$ws.Cells.Item(40, 5).Value = '  +2.09%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.9054'
$ws.Cells.Item(41, 5).Value = '  +0.46%  '

$ws.Cells.Item(42, 5).Value = '  +0.11%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '2.002.86'
$ws.Cells.Item(43, 5).Value = '  +0.06%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '101.66'
$ws.Cells.Item(44, 5).Value = '  -0.10%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '65.84'
$ws.Cells.Item(45, 5).Value = '  +0.04%  '

$ws.Cells.Item(46, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.00000000120'
$ws.Cells.Item(46, 5).Value = '  +1.21%  '

$ws.Cells.Item(47, 2).Value = 'Aptos'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '7.152'
$ws.Cells.Item(47, 5).Value = '  +0.79%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.4008'
$ws.Cells.Item(48, 5).Value = '  +0.28%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.1157'
$ws.Cells.Item(49, 5).Value = '  +3.26%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '9.010'
$ws.Cells.Item(50, 5).Value = '  +1.33%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '1.672'
$ws.Cells.Item(51, 5).Value = '  +0.06%  '
